$wb = $excel.ActiveWorkbook

# Set the PDiCECpDoC!B2 cell to a hardcoded value instead of the formula
# referencing the "Texas Notes" sheet that is about to be removed.
$pd = $wb.Worksheets.Item("PDiCECpDoC")
$pd.Range("B2").Value = 0.13

# Delete the "Texas Notes" sheet entirely.
$excel.DisplayAlerts = $false
$tn = $wb.Worksheets.Item("Texas Notes")
$tn.Delete()
$excel.DisplayAlerts = $true
